$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.949.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.429.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.20%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.515"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.17%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("E12").Value = "  -2.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.801.23"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.416.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.832"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.765.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("E25").Value = "  +2.39%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +16.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.22"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.86%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0765"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "127.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.940.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.78"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +17.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.57%  "
